# Actualización automática 2025-08-01 08:30:08
#
# Applies the monthly rollover update to the "HIDALGO HIDALGO PEDRO GUSTAVO"
# workbook:
#   - "VENTAS POR GRUPO": zero out the current-period (last month column
#     group) figures that have now rolled off the 4-month rolling window,
#     and refresh the "N de 20" completion counters in the totals row.
#   - "VENTA MENSUAL": shift the 4 rolling month columns (the header labels
#     move from abril/mayo/junio/julio to mayo/junio/julio/agosto) and
#     refresh the figures + column widths for the new period.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": zero out the cells that rolled off the window
# ---------------------------------------------------------------------
$ws1.Range("E5").Value  = 0
$ws1.Range("M5").Value  = 0

$ws1.Range("D6").Value  = 0

$ws1.Range("D7").Value  = 0
$ws1.Range("L7").Value  = 0
$ws1.Range("M7").Value  = 0
$ws1.Range("O7").Value  = 0

$ws1.Range("M8").Value  = 0
$ws1.Range("Q8").Value  = 0

$ws1.Range("H10").Value = 0
$ws1.Range("I10").Value = 0
$ws1.Range("M10").Value = 0
$ws1.Range("P10").Value = 0

$ws1.Range("M12").Value = 0

$ws1.Range("I13").Value = 0
$ws1.Range("M13").Value = 0

$ws1.Range("H15").Value = 0
$ws1.Range("I15").Value = 0
$ws1.Range("M15").Value = 0
$ws1.Range("Q15").Value = 0

$ws1.Range("M16").Value = 0

$ws1.Range("M18").Value = 0

$ws1.Range("M20").Value = 0

$ws1.Range("D21").Value = 0
$ws1.Range("H21").Value = 0
$ws1.Range("I21").Value = 0
$ws1.Range("M21").Value = 0

# Totals row: refresh the "N de 20" completion counters that moved to 0
$ws1.Range("D22").Value = "0 de 20"
$ws1.Range("E22").Value = "0 de 20"
$ws1.Range("H22").Value = "0 de 20"
$ws1.Range("I22").Value = "0 de 20"
$ws1.Range("L22").Value = "0 de 20"
$ws1.Range("M22").Value = "0 de 20"
$ws1.Range("O22").Value = "0 de 20"
$ws1.Range("P22").Value = "0 de 20"
$ws1.Range("Q22").Value = "0 de 20"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": roll the 4-month window forward by one month
# ---------------------------------------------------------------------

# Column widths for C:F (raw OOXML width = ColumnWidth + 0.83 in this
# engine's unit system, matching real Excel's Normal-style-font offset)
$ws2.Columns.Item(3).ColumnWidth = 13.17
$ws2.Columns.Item(4).ColumnWidth = 12.17
$ws2.Columns.Item(5).ColumnWidth = 13.17
$ws2.Columns.Item(6).ColumnWidth = 11.17

# Month headers
$ws2.Range("C1").Value = "mayo"
$ws2.Range("D1").Value = "junio"
$ws2.Range("E1").Value = "julio"
$ws2.Range("F1").Value = "agosto"

# Row 5 - CARRION CARRION LESLY ANABE
$ws2.Range("C5").Value  = 1914.36
$ws2.Range("D5").Value  = 12.1
$ws2.Range("E5").Value  = 9064.610000000001
$ws2.Range("F5").Value  = 0

# Row 6 - CHASIQUIZA CAMPAÑA JOSE LUIS
$ws2.Range("C6").Value  = 6839.4
$ws2.Range("D6").Value  = 111.02
$ws2.Range("E6").Value  = 7193.01
$ws2.Range("F6").Value  = 0

# Row 7 - CHONTASI SIMBAÑA SILVIA JANETH
$ws2.Range("C7").Value  = 595.8
$ws2.Range("D7").Value  = 177.41
$ws2.Range("E7").Value  = 2898.45
$ws2.Range("F7").Value  = 0

# Row 8 - DECORHOME S.C.C.
$ws2.Range("C8").Value  = 7010.57
$ws2.Range("D8").Value  = 211.12
$ws2.Range("E8").Value  = 2279.62
$ws2.Range("F8").Value  = 0

# Row 10 - JARAMILLO CARVAJAL NICOLAS ESTEBAN
$ws2.Range("C10").Value = 11565.71
$ws2.Range("D10").Value = 4711.21
$ws2.Range("E10").Value = 12314.93
$ws2.Range("F10").Value = 0

# Row 12 - MEGAMAFERS S.A.
$ws2.Range("C12").Value = 1565.15
$ws2.Range("D12").Value = 4381.39
$ws2.Range("E12").Value = 3399.19
$ws2.Range("F12").Value = 0

# Row 13 - MUÑOZ LOZA ROMMEL SEBASTIAN
$ws2.Range("C13").Value = 9034.549999999999
$ws2.Range("D13").Value = 2291.68
$ws2.Range("E13").Value = 5370.47
$ws2.Range("F13").Value = 0

# Row 15 - OÑATE PEREZ MERCY YOLANDA
$ws2.Range("C15").Value = 288.17
$ws2.Range("D15").Value = 431.64
$ws2.Range("E15").Value = 1319.85
$ws2.Range("F15").Value = 0

# Row 16 - PADILLA MIER BERTHA MARIETA
$ws2.Range("C16").Value = 0
$ws2.Range("D16").Value = 45.91
$ws2.Range("E16").Value = 12316.05
$ws2.Range("F16").Value = 0

# Row 18 - SARZOSA UNDA JOSE DOMINGO
$ws2.Range("C18").Value = 290.56
$ws2.Range("D18").Value = 1935.16
$ws2.Range("E18").Value = 1702.17
$ws2.Range("F18").Value = 0

# Row 20 - TRUJILLO TORRES VINICIO RUBEN
$ws2.Range("E20").Value = 2925.73
$ws2.Range("F20").Value = 0

# Row 21 - TULCAN NARVAEZ EDITH MARITZA
$ws2.Range("C21").Value = 15564.83
$ws2.Range("D21").Value = 1964.99
$ws2.Range("E21").Value = 5388.82
$ws2.Range("F21").Value = 0

# Row 22 - TOTAL
$ws2.Range("C22").Value = 54669.1
$ws2.Range("D22").Value = 16273.63
$ws2.Range("E22").Value = 66172.89999999999
$ws2.Range("F22").Value = 0
